$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting from the last existing data row (63) onto the two new rows (64, 65) ---
$ws.Range("A63:V63").Copy() | Out-Null
$ws.Range("A64:V64").PasteSpecial(-4122) | Out-Null
$ws.Range("A63:V63").Copy() | Out-Null
$ws.Range("A65:V65").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 64: HNK Gorica 3 x 0 Osijek ---
$ws.Cells.Item(64, 1).Value = 63
$ws.Cells.Item(64, 2).Value = "croatia"
$ws.Cells.Item(64, 3).Value = "hnl"
$ws.Cells.Item(64, 4).Value = "2023-2024"
$ws.Cells.Item(64, 5).Value = 45234.66666666666
$ws.Cells.Item(64, 6).Value = "Gorica"
$ws.Cells.Item(64, 7).Value = 3
$ws.Cells.Item(64, 8).Value = "Osijek"
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 2.79
$ws.Cells.Item(64, 11).Value = "29/10/2023 15:12"
$ws.Cells.Item(64, 12).Value = 2.72
$ws.Cells.Item(64, 13).Value = "04/11/2023 15:59"
$ws.Cells.Item(64, 14).Value = 3.27
$ws.Cells.Item(64, 15).Value = "29/10/2023 15:12"
$ws.Cells.Item(64, 16).Value = 3.06
$ws.Cells.Item(64, 17).Value = "04/11/2023 15:59"
$ws.Cells.Item(64, 18).Value = 2.57
$ws.Cells.Item(64, 19).Value = "29/10/2023 15:12"
$ws.Cells.Item(64, 20).Value = 2.86
$ws.Cells.Item(64, 21).Value = "04/11/2023 15:59"
$ws.Cells.Item(64, 22).Value = "https://www.betexplorer.com/football/croatia/hnl/hnk-gorica-osijek/xbeMXuKg/"

# --- Row 65: Hajduk Split 1 x 0 Rudes ---
$ws.Cells.Item(65, 1).Value = 64
$ws.Cells.Item(65, 2).Value = "croatia"
$ws.Cells.Item(65, 3).Value = "hnl"
$ws.Cells.Item(65, 4).Value = "2023-2024"
$ws.Cells.Item(65, 5).Value = 45234.75694444445
$ws.Cells.Item(65, 6).Value = "Hajduk Split"
$ws.Cells.Item(65, 7).Value = 1
$ws.Cells.Item(65, 8).Value = "Rudes"
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 1.19
$ws.Cells.Item(65, 11).Value = "28/10/2023 19:13"
$ws.Cells.Item(65, 12).Value = 1.17
$ws.Cells.Item(65, 13).Value = "04/11/2023 18:06"
$ws.Cells.Item(65, 14).Value = 6.93
$ws.Cells.Item(65, 15).Value = "28/10/2023 19:13"
$ws.Cells.Item(65, 16).Value = 7.03
$ws.Cells.Item(65, 17).Value = "04/11/2023 18:06"
$ws.Cells.Item(65, 18).Value = 12.66
$ws.Cells.Item(65, 19).Value = "28/10/2023 19:13"
$ws.Cells.Item(65, 20).Value = 16.62
$ws.Cells.Item(65, 21).Value = "04/11/2023 18:06"
$ws.Cells.Item(65, 22).Value = "https://www.betexplorer.com/football/croatia/hnl/hajduk-split-rudes/YPlVVJl6/"
